$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "27.995.73"
Set-TextCell "E2" "  -0.26%  "
Set-TextCell "D3" "1.857.97"
Set-TextCell "E3" "  -0.89%  "
Set-TextCell "E4" "  +0.10%  "
Set-TextCell "D5" "312.44"
Set-TextCell "E5" "  -0.37%  "
Set-TextCell "E6" "  +0.04%  "
Set-TextCell "D7" "0.5142"
Set-TextCell "E7" "  +1.50%  "
Set-TextCell "D8" "0.3838"
Set-TextCell "E8" "  -0.06%  "
Set-TextCell "D9" "0.08218"
Set-TextCell "E9" "  -8.30%  "
Set-TextCell "D10" "1.109"
Set-TextCell "E10" "  -1.04%  "
Set-TextCell "D11" "41.49"
Set-TextCell "E11" "  -0.28%  "
Set-TextCell "D12" "6.186"
Set-TextCell "E12" "  -2.24%  "
Set-TextCell "E13" "  -0.69%  "
Set-TextCell "D14" "1.866.39"
Set-TextCell "E14" "  -0.42%  "
Set-TextCell "D15" "7.255"
Set-TextCell "E15" "  +0.97%  "
Set-TextCell "E16" "  +0.04%  "
Set-TextCell "D17" "0.00001097"
Set-TextCell "E17" "  -0.79%  "
Set-TextCell "D18" "90.51"
Set-TextCell "E18" "  -0.55%  "
Set-TextCell "D19" "0.06655"
Set-TextCell "E19" "  +0.86%  "
Set-TextCell "D20" "17.68"
Set-TextCell "E20" "  -2.32%  "
Set-TextCell "E21" "  +0.03%  "
Set-TextCell "D22" "6.001"
Set-TextCell "E22" "  -1.72%  "
Set-TextCell "D23" "28.019.56"
Set-TextCell "E23" "  -0.23%  "
Set-TextCell "E24" "  -2.97%  "
Set-TextCell "D25" "2.245"
Set-TextCell "E25" "  -1.63%  "
Set-TextCell "D26" "2.075.12"
Set-TextCell "E26" "  -0.70%  "
Set-TextCell "D27" "2.508"
Set-TextCell "E27" "  -0.92%  "
Set-TextCell "D28" "157.97"
Set-TextCell "E28" "  +0.61%  "
Set-TextCell "D29" "20.45"
Set-TextCell "E29" "  -1.42%  "
Set-TextCell "D30" "124.49"
Set-TextCell "E30" "  -1.62%  "
Set-TextCell "D31" "0.1065"
Set-TextCell "E31" "  +1.35%  "
Set-TextCell "E32" "  -2.80%  "
Set-TextCell "D33" "5.984"
Set-TextCell "E33" "  +6.83%  "
Set-TextCell "D34" "3.599"
Set-TextCell "E34" "  -0.17%  "
Set-TextCell "D35" "9.346"
Set-TextCell "E35" "  -2.90%  "
Set-TextCell "D36" "0.02415"
Set-TextCell "E36" "  -0.23%  "
Set-TextCell "D37" "0.06492"
Set-TextCell "E37" "  -1.24%  "
Set-TextCell "E38" "  -0.43%  "
Set-TextCell "D39" "0.6536"
Set-TextCell "E39" "  +2.45%  "
Set-TextCell "D40" "1.194"
Set-TextCell "E40" "  -0.94%  "
Set-TextCell "D41" "5.011"
Set-TextCell "E41" "  +2.11%  "
Set-TextCell "D42" "1.218"
Set-TextCell "E42" "  -3.84%  "
Set-TextCell "E43" "  -2.48%  "
Set-TextCell "D44" "0.6138"
Set-TextCell "E44" "  +2.16%  "
Set-TextCell "D45" "12.97"
Set-TextCell "E45" "  -1.89%  "
Set-TextCell "D46" "1.281"
Set-TextCell "E46" "  +0.36%  "
Set-TextCell "D47" "3.665"
Set-TextCell "E47" "  -0.17%  "
Set-TextCell "D48" "2.008"
Set-TextCell "E48" "  +0.95%  "
Set-TextCell "D50" "120.34"
Set-TextCell "E50" "  -0.78%  "
Set-TextCell "D51" "78.39"
Set-TextCell "E51" "  -1.55%  "
